# ajout des figures e- au rapport dosi relative
#
# - Renomme l'onglet "CR_09MeV_6_DSP100_CC13" en "CR_09MeV_06_DSP100_CC13"
#   (harmonisation du format du nombre de cm, cf. "_06_" / "_09_" / "_15_" ...)
# - Active cet onglet (il devient l'onglet selectionne / actif du classeur)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CR_09MeV_6_DSP100_CC13")
$ws.Name = "CR_09MeV_06_DSP100_CC13"

# Rend cet onglet actif / selectionne (met a jour activeTab du classeur et
# tabSelected de la feuille, et retire tabSelected de l'ancien onglet actif).
$ws.Activate()
